$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Clear the "Results" column (J2:J12), which previously held "Pass"/"Fail"
# values. This removes the now-unused "Pass"/"Fail" shared strings.
$rng = $ws.Range("J2:J12")
$rng.ClearContents()

# Reflect the resulting selection as seen in the saved workbook.
$ws.Activate()
$rng.Select()
